$d = $word.ActiveDocument

# 1. "Programa resumido" paragraph: split after "orientador," with a manual line break
$d.Content.Find.Execute(
    "Desenvolvimento do trabalho de conclusão de curso, sob orientação de um professor orientador,o qual deve constituir-se num projeto de tema específico relacionado às atribuições da profissão.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Desenvolvimento do trabalho de conclusão de curso, sob orientação de um professor orientador,^lo qual deve constituir-se num projeto de tema específico relacionado às atribuições da profissão.",
    2)

# 2. "Programa" paragraph: split into four runs at the missing-space boundaries
$d.Content.Find.Execute(
    "Elaboração de uma monografia ou de relatório técnico que apresente: (1) o tema e suaimportância, (2) os objetivos, (3) a revisão bibliográfica, (4) a metodologia científica (5) odesenvolvimento do projeto, (6) a análise e discussão dos resultados, (7) as conclusões e (8)referências bibliográficas",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Elaboração de uma monografia ou de relatório técnico que apresente: (1) o tema e sua^limportância, (2) os objetivos, (3) a revisão bibliográfica, (4) a metodologia científica (5) o^ldesenvolvimento do projeto, (6) a análise e discussão dos resultados, (7) as conclusões e (8)^lreferências bibliográficas",
    2)

# 3. "Método" run: split into three runs at the missing-space boundaries
$d.Content.Find.Execute(
    "Reuniões periódicas com o orientador e realização do trabalho conforme orientação eapresentação de uma monografia final, conforme norma do Departamento de Engenharia Químicae Produção.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Reuniões periódicas com o orientador e realização do trabalho conforme orientação e^lapresentação de uma monografia final, conforme norma do Departamento de Engenharia Química^le Produção.",
    2)

# 4. "Critério" run: split into two runs at the missing-space boundary
$d.Content.Find.Execute(
    "Avaliação perante uma banca examinadora composta por 3 (três) membros, conforme norma doDepartamento de Engenharia Química e Produção.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Avaliação perante uma banca examinadora composta por 3 (três) membros, conforme norma do^lDepartamento de Engenharia Química e Produção.",
    2)
